$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Edit 1: "Other Members Present: ... Gerjan Haxhia" -> "... Gerjan Haxhija"
#         (fix the misspelled surname by inserting a "j" before the final "a")
# -----------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Gerjan Haxhia", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    # $rng now covers "Gerjan Haxhia"; collapse to its end, then step back
    # one character so the insertion point sits right before the final "a".
    $rng.Collapse(0)
    $rng.MoveStart(1, -1) | Out-Null
    $rng.Collapse(1)
    $rng.InsertBefore("j")
    # Nudge formatting off/on so the newly typed "j" keeps its own run,
    # matching how Word splits a run when you type in the middle of it.
    $rng.Bold = 0
    $rng.Bold = 1
}

# -----------------------------------------------------------------------
# Edit 2: the "goal for this week" paragraph was split across three runs
#         ("...focusing on " / "preparing for " / "Activity Diagrams...");
#         re-merge them into a single contiguous run.
# -----------------------------------------------------------------------
$goalText = "The goal for this week is to continue working on the requirements " + `
            "specification document, focusing on preparing for Activity Diagrams, " + `
            "State Diagrams, Sequence Diagrams, and Interaction Diagrams. The " + `
            "following task assignments have been made: "
$d.Content.Find.Execute($goalText, $false, $false, $false, $false, $false, `
                         $true, 1, $false, $goalText, 2) | Out-Null

# -----------------------------------------------------------------------
# Edit 3: the blank paragraph right after "Xhoni: Use Case Diagrams" had an
#         empty run (<w:rPr/>) with no explicit formatting; give it the
#         same Arial/non-bold/11pt formatting used by the other body runs.
# -----------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -match "Xhoni: Use Case Diagrams") {
        $blank = $paras.Item($i + 1).Range
        $blank.Font.Name = "Arial"
        $blank.Font.Bold = 0
        $blank.Font.Size = 11
        break
    }
}
